$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header columns V1:Y1 (bold/bordered header style carried automatically by column) ---
$ws.Range("V1").Value2 = "Posesión Local (%)"
$ws.Range("W1").Value2 = "Posesión Visita (%)"
$ws.Range("X1").Value2 = "fuente_tiempos"
$ws.Range("Y1").Value2 = "estado_datos"
$ws.Range("U1").Copy() | Out-Null
$ws.Range("V1:Y1").PasteSpecial(-4122) | Out-Null
$ws.Range("V1").Value2 = "Posesión Local (%)"
$ws.Range("W1").Value2 = "Posesión Visita (%)"
$ws.Range("X1").Value2 = "fuente_tiempos"
$ws.Range("Y1").Value2 = "estado_datos"

# --- 2. Scattered data corrections on existing rows (goles 1T/2T swaps, cleared cells) ---
$ws.Cells.Item(155,13).Value2 = 1
$ws.Cells.Item(155,15).Value2 = 0
$ws.Cells.Item(156,13).Value2 = 1
$ws.Cells.Item(156,15).Value2 = 0
$ws.Cells.Item(159,13).Value2 = 1
$ws.Cells.Item(159,15).Value2 = 1
$ws.Cells.Item(160,13).Value2 = 4
$ws.Cells.Item(160,14).Value2 = 1
$ws.Cells.Item(160,15).Value2 = 1
$ws.Cells.Item(160,16).Value2 = 0
$ws.Cells.Item(161,13).Value2 = 1
$ws.Cells.Item(161,15).Value2 = 0
$ws.Cells.Item(169,13).Value2 = 1
$ws.Cells.Item(169,14).Value2 = 1
$ws.Cells.Item(169,15).Value2 = 1
$ws.Cells.Item(169,16).Value2 = 0
$ws.Cells.Item(172,13).Value2 = 2
$ws.Cells.Item(172,15).Value2 = 1
$ws.Cells.Item(174,13).Value2 = 1
$ws.Cells.Item(174,15).Value2 = 2
$ws.Cells.Item(175,13).Value2 = 2
$ws.Cells.Item(175,14).Value2 = 1
$ws.Cells.Item(175,15).Value2 = 1
$ws.Cells.Item(175,16).Value2 = 0
$ws.Cells.Item(176,14).Value2 = 1
$ws.Cells.Item(176,16).Value2 = 1
$ws.Cells.Item(179,13).Value2 = 2
$ws.Cells.Item(179,15).Value2 = 1
$ws.Cells.Item(181,13).Value2 = 2
$ws.Cells.Item(181,15).Value2 = 1
$ws.Cells.Item(182,13).Value2 = 3
$ws.Cells.Item(182,15).Value2 = 2
$ws.Cells.Item(184,13).Value2 = 1
$ws.Cells.Item(184,15).Value2 = 0
$ws.Cells.Item(185,14).Value2 = 1
$ws.Cells.Item(185,16).Value2 = 1
$ws.Cells.Item(187,14).Value2 = 1
$ws.Cells.Item(187,16).Value2 = 1
$ws.Cells.Item(191,13).Value2 = 2
$ws.Cells.Item(191,15).Value2 = 3
$ws.Cells.Item(192,14).Value2 = 1
$ws.Cells.Item(192,16).Value2 = 2
$ws.Cells.Item(193,13).Value2 = 1
$ws.Cells.Item(193,15).Value2 = 1
$ws.Cells.Item(194,13).Value2 = 1
$ws.Cells.Item(194,14).Value2 = 1
$ws.Cells.Item(194,15).Value2 = 0
$ws.Cells.Item(194,16).Value2 = 1
$ws.Cells.Item(195,13).Value2 = 2
$ws.Cells.Item(195,15).Value2 = 1
$ws.Cells.Item(196,13).Value2 = 2
$ws.Cells.Item(196,14).Value2 = 2
$ws.Cells.Item(196,15).Value2 = 1
$ws.Cells.Item(196,16).Value2 = 0
$ws.Cells.Item(198,13).Value2 = 1
$ws.Cells.Item(198,14).Value2 = 1
$ws.Cells.Item(198,15).Value2 = 1
$ws.Cells.Item(198,16).Value2 = 0
$ws.Cells.Item(199,13).Value2 = 2
$ws.Cells.Item(199,14).Value2 = 1
$ws.Cells.Item(199,15).Value2 = 1
$ws.Cells.Item(199,16).Value2 = 0
$ws.Cells.Item(200,13).Value2 = 1
$ws.Cells.Item(200,14).Value2 = 1
$ws.Cells.Item(200,15).Value2 = 0
$ws.Cells.Item(200,16).Value2 = 1
$ws.Cells.Item(201,13).Value2 = 1
$ws.Cells.Item(201,15).Value2 = 0
$ws.Cells.Item(202,13).Value2 = 1
$ws.Cells.Item(202,15).Value2 = 2
$ws.Cells.Item(203,14).Value2 = 1
$ws.Cells.Item(203,16).Value2 = 1
$ws.Cells.Item(204,13).Value2 = 1
$ws.Cells.Item(204,15).Value2 = 0
$ws.Cells.Item(207,13).Value2 = 2
$ws.Cells.Item(207,15).Value2 = 0

# --- 3. Clear cells that no longer have data (now blank) ---
$ws.Cells.Item(162,7).Value2 = ""
$ws.Cells.Item(162,8).Value2 = ""
$ws.Cells.Item(162,9).Value2 = ""
$ws.Cells.Item(162,10).Value2 = ""
$ws.Cells.Item(162,17).Value2 = ""
$ws.Cells.Item(162,18).Value2 = ""
$ws.Cells.Item(173,9).Value2 = ""
$ws.Cells.Item(173,10).Value2 = ""
$ws.Cells.Item(203,9).Value2 = ""
$ws.Cells.Item(203,10).Value2 = ""

# --- 4. Append new fixture rows 208-216 ---
# row 208: Comerciantes Unidos vs Cienciano
$ws.Cells.Item(208,1).Value2 = "'2025-08-08"
$ws.Cells.Item(208,2).Value2 = "Comerciantes Unidos"
$ws.Cells.Item(208,3).Value2 = "Cienciano"
$ws.Cells.Item(208,4).Value2 = 1
$ws.Cells.Item(208,5).Value2 = 1
$ws.Cells.Item(208,6).Value2 = 1405749
$ws.Cells.Item(208,7).Value2 = 2
$ws.Cells.Item(208,8).Value2 = 7
$ws.Cells.Item(208,9).Value2 = 2
$ws.Cells.Item(208,10).Value2 = 2
$ws.Cells.Item(208,11).Value2 = 0
$ws.Cells.Item(208,12).Value2 = 0
$ws.Cells.Item(208,13).Value2 = 0
$ws.Cells.Item(208,14).Value2 = 0
$ws.Cells.Item(208,15).Value2 = 1
$ws.Cells.Item(208,16).Value2 = 1
$ws.Cells.Item(208,17).Value2 = 36
$ws.Cells.Item(208,18).Value2 = 64
$ws.Cells.Item(208,19).Value2 = "E"

# row 209: Ayacucho FC vs Alianza Lima
$ws.Cells.Item(209,1).Value2 = "'2025-08-09"
$ws.Cells.Item(209,2).Value2 = "Ayacucho FC"
$ws.Cells.Item(209,3).Value2 = "Alianza Lima"
$ws.Cells.Item(209,4).Value2 = 0
$ws.Cells.Item(209,5).Value2 = 1
$ws.Cells.Item(209,6).Value2 = 1405750
$ws.Cells.Item(209,7).Value2 = 0
$ws.Cells.Item(209,8).Value2 = 1
$ws.Cells.Item(209,9).Value2 = ""
$ws.Cells.Item(209,10).Value2 = ""
$ws.Cells.Item(209,11).Value2 = 0
$ws.Cells.Item(209,12).Value2 = 0
$ws.Cells.Item(209,13).Value2 = 0
$ws.Cells.Item(209,14).Value2 = 0
$ws.Cells.Item(209,15).Value2 = 0
$ws.Cells.Item(209,16).Value2 = 1
$ws.Cells.Item(209,17).Value2 = 38
$ws.Cells.Item(209,18).Value2 = 62
$ws.Cells.Item(209,19).Value2 = "V"

# row 210: ADT vs Sport Huancayo
$ws.Cells.Item(210,1).Value2 = "'2025-08-09"
$ws.Cells.Item(210,2).Value2 = "ADT"
$ws.Cells.Item(210,3).Value2 = "Sport Huancayo"
$ws.Cells.Item(210,4).Value2 = 1
$ws.Cells.Item(210,5).Value2 = 0
$ws.Cells.Item(210,6).Value2 = 1429385
$ws.Cells.Item(210,7).Value2 = 5
$ws.Cells.Item(210,8).Value2 = 3
$ws.Cells.Item(210,9).Value2 = 2
$ws.Cells.Item(210,10).Value2 = 3
$ws.Cells.Item(210,11).Value2 = 1
$ws.Cells.Item(210,12).Value2 = 0
$ws.Cells.Item(210,13).Value2 = 0
$ws.Cells.Item(210,14).Value2 = 0
$ws.Cells.Item(210,15).Value2 = 1
$ws.Cells.Item(210,16).Value2 = 0
$ws.Cells.Item(210,17).Value2 = 53
$ws.Cells.Item(210,18).Value2 = 47
$ws.Cells.Item(210,19).Value2 = "L"

# row 211: Universitario vs Sport Boys
$ws.Cells.Item(211,1).Value2 = "'2025-08-09"
$ws.Cells.Item(211,2).Value2 = "Universitario"
$ws.Cells.Item(211,3).Value2 = "Sport Boys"
$ws.Cells.Item(211,4).Value2 = 1
$ws.Cells.Item(211,5).Value2 = 0
$ws.Cells.Item(211,6).Value2 = 1405751
$ws.Cells.Item(211,7).Value2 = 10
$ws.Cells.Item(211,8).Value2 = 2
$ws.Cells.Item(211,9).Value2 = 1
$ws.Cells.Item(211,10).Value2 = 1
$ws.Cells.Item(211,11).Value2 = 0
$ws.Cells.Item(211,12).Value2 = 0
$ws.Cells.Item(211,13).Value2 = 1
$ws.Cells.Item(211,14).Value2 = 0
$ws.Cells.Item(211,15).Value2 = 0
$ws.Cells.Item(211,16).Value2 = 0
$ws.Cells.Item(211,17).Value2 = 56
$ws.Cells.Item(211,18).Value2 = 44
$ws.Cells.Item(211,19).Value2 = "L"

# row 212: Sporting Cristal vs FBC Melgar
$ws.Cells.Item(212,1).Value2 = "'2025-08-10"
$ws.Cells.Item(212,2).Value2 = "Sporting Cristal"
$ws.Cells.Item(212,3).Value2 = "FBC Melgar"
$ws.Cells.Item(212,4).Value2 = 1
$ws.Cells.Item(212,5).Value2 = 0
$ws.Cells.Item(212,6).Value2 = 1405752
$ws.Cells.Item(212,7).Value2 = 0
$ws.Cells.Item(212,8).Value2 = 2
$ws.Cells.Item(212,9).Value2 = 0
$ws.Cells.Item(212,10).Value2 = 4
$ws.Cells.Item(212,11).Value2 = 1
$ws.Cells.Item(212,12).Value2 = 1
$ws.Cells.Item(212,13).Value2 = 1
$ws.Cells.Item(212,14).Value2 = 0
$ws.Cells.Item(212,15).Value2 = 0
$ws.Cells.Item(212,16).Value2 = 0
$ws.Cells.Item(212,17).Value2 = 63
$ws.Cells.Item(212,18).Value2 = 37
$ws.Cells.Item(212,19).Value2 = "L"

# row 213: Cultural Santa Rosa vs Atletico Grau
$ws.Cells.Item(213,1).Value2 = "'2025-08-10"
$ws.Cells.Item(213,2).Value2 = "Cultural Santa Rosa"
$ws.Cells.Item(213,3).Value2 = "Atletico Grau"
$ws.Cells.Item(213,4).Value2 = 3
$ws.Cells.Item(213,5).Value2 = 2
$ws.Cells.Item(213,6).Value2 = 1410135
$ws.Cells.Item(213,7).Value2 = 1
$ws.Cells.Item(213,8).Value2 = 3
$ws.Cells.Item(213,9).Value2 = 1
$ws.Cells.Item(213,10).Value2 = 1
$ws.Cells.Item(213,11).Value2 = 0
$ws.Cells.Item(213,12).Value2 = 1
$ws.Cells.Item(213,13).Value2 = 0
$ws.Cells.Item(213,14).Value2 = 1
$ws.Cells.Item(213,15).Value2 = 3
$ws.Cells.Item(213,16).Value2 = 1
$ws.Cells.Item(213,17).Value2 = 43
$ws.Cells.Item(213,18).Value2 = 57
$ws.Cells.Item(213,19).Value2 = "L"

# row 214: Juan Pablo II College vs Deportivo Binacional
$ws.Cells.Item(214,1).Value2 = "'2025-08-10"
$ws.Cells.Item(214,2).Value2 = "Juan Pablo II College"
$ws.Cells.Item(214,3).Value2 = "Deportivo Binacional"
$ws.Cells.Item(214,4).Value2 = 0
$ws.Cells.Item(214,5).Value2 = 0
$ws.Cells.Item(214,6).Value2 = 1410136
$ws.Cells.Item(214,7).Value2 = 0
$ws.Cells.Item(214,8).Value2 = 5
$ws.Cells.Item(214,9).Value2 = 2
$ws.Cells.Item(214,10).Value2 = 1
$ws.Cells.Item(214,11).Value2 = 0
$ws.Cells.Item(214,12).Value2 = 0
$ws.Cells.Item(214,13).Value2 = 0
$ws.Cells.Item(214,14).Value2 = 0
$ws.Cells.Item(214,15).Value2 = 0
$ws.Cells.Item(214,16).Value2 = 0
$ws.Cells.Item(214,17).Value2 = 51
$ws.Cells.Item(214,18).Value2 = 49
$ws.Cells.Item(214,19).Value2 = "E"

# row 215: Alianza Atletico vs UTC
$ws.Cells.Item(215,1).Value2 = "'2025-08-11"
$ws.Cells.Item(215,2).Value2 = "Alianza Atletico"
$ws.Cells.Item(215,3).Value2 = "UTC"
$ws.Cells.Item(215,4).Value2 = 2
$ws.Cells.Item(215,5).Value2 = 0
$ws.Cells.Item(215,6).Value2 = 1405753
$ws.Cells.Item(215,7).Value2 = 4
$ws.Cells.Item(215,8).Value2 = 3
$ws.Cells.Item(215,9).Value2 = 2
$ws.Cells.Item(215,10).Value2 = 1
$ws.Cells.Item(215,11).Value2 = 1
$ws.Cells.Item(215,12).Value2 = 0
$ws.Cells.Item(215,13).Value2 = 2
$ws.Cells.Item(215,14).Value2 = 0
$ws.Cells.Item(215,15).Value2 = 0
$ws.Cells.Item(215,16).Value2 = 0
$ws.Cells.Item(215,17).Value2 = 48
$ws.Cells.Item(215,18).Value2 = 52
$ws.Cells.Item(215,19).Value2 = "L"

# row 216: Deportivo Garcilaso vs Alianza Universidad
$ws.Cells.Item(216,1).Value2 = "'2025-08-12"
$ws.Cells.Item(216,2).Value2 = "Deportivo Garcilaso"
$ws.Cells.Item(216,3).Value2 = "Alianza Universidad"
$ws.Cells.Item(216,4).Value2 = 1
$ws.Cells.Item(216,5).Value2 = 1
$ws.Cells.Item(216,6).Value2 = 1405754
$ws.Cells.Item(216,7).Value2 = 6
$ws.Cells.Item(216,8).Value2 = 2
$ws.Cells.Item(216,9).Value2 = 3
$ws.Cells.Item(216,10).Value2 = 4
$ws.Cells.Item(216,11).Value2 = 0
$ws.Cells.Item(216,12).Value2 = 0
$ws.Cells.Item(216,13).Value2 = 0
$ws.Cells.Item(216,14).Value2 = 0
$ws.Cells.Item(216,15).Value2 = 1
$ws.Cells.Item(216,16).Value2 = 1
$ws.Cells.Item(216,17).Value2 = 48
$ws.Cells.Item(216,18).Value2 = 52
$ws.Cells.Item(216,19).Value2 = "E"

Write-Host "Edit complete."
